# Reran all processing steps after manually fixing wrong recordings date
# information based on app logs -- this updates the statistics results
# (normality, equal_var, mixed_anova, pairwise_ttests) with the recomputed
# values.

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell while forcing it to stay a *text*
# cell (the pairwise_ttests "BF10" column stores numbers formatted as
# strings, e.g. "2.32e+21" -- a plain numeric assignment would silently
# convert it to a real number). We temporarily mark the cell as Text,
# assign the string, then restore the default "Normal" style so no
# stray number-format sticks around on the cell.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------
# Sheet "normality"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("normality")
$ws.Range("C3").Value = 0.9566
$ws.Range("D3").Value = 0.0005
$ws.Range("C4").Value = 0.954
$ws.Range("D4").Value = 0.07729999999999999
$ws.Range("C5").Value = 0.9668
$ws.Range("D5").Value = 0.0036
$ws.Range("C6").Value = 0.9762999999999999
$ws.Range("D6").Value = 0.4912
$ws.Range("D7").Value = 0.0059
$ws.Range("C8").Value = 0.9643
$ws.Range("D8").Value = 0.1881
$ws.Range("C9").Value = 0.9608
$ws.Range("D9").Value = 0.0011
$ws.Range("C10").Value = 0.9658
$ws.Range("D10").Value = 0.2143
$ws.Range("C11").Value = 0.93
$ws.Range("C12").Value = 0.9478
$ws.Range("D12").Value = 0.0456

# ---------------------------------------------------------------
# Sheet "equal_var"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("equal_var")
$ws.Range("C3").Value = 4.7792
$ws.Range("D3").Value = 0.0302
$ws.Range("C4").Value = 3.9715
$ws.Range("D4").Value = 0.0479
$ws.Range("C5").Value = 2.9453
$ws.Range("D5").Value = 0.08799999999999999
$ws.Range("C6").Value = 5.3489
$ws.Range("D6").Value = 0.022
$ws.Range("C7").Value = 2.5608
$ws.Range("D7").Value = 0.1114

# ---------------------------------------------------------------
# Sheet "mixed_anova"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("mixed_anova")
$ws.Range("C3").Value = 2.4713
$ws.Range("E3").Value = 167
$ws.Range("F3").Value = 2.4713
$ws.Range("G3").Value = 0.0257
$ws.Range("H3").Value = 0.8729

$ws.Range("C4").Value = 4063.2874
$ws.Range("E4").Value = 668
$ws.Range("F4").Value = 1015.8218
$ws.Range("G4").Value = 67.9883
$ws.Range("J4").Value = 0.2893
$ws.Range("K4").Value = 0.4913
$ws.Range("M4").Value = 0.1004

$ws.Range("C5").Value = 15.8148
$ws.Range("E5").Value = 668
$ws.Range("F5").Value = 3.9537
$ws.Range("G5").Value = 0.2646
$ws.Range("H5").Value = 0.9006999999999999
$ws.Range("J5").Value = 0.0016

# ---------------------------------------------------------------
# Sheet "pairwise_ttests"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("pairwise_ttests")

$ws.Range("H3").Value = -12.0498
$ws.Range("I3").Value = 168
Set-TextValue $ws.Range("L3") "2.32e+21"
$ws.Range("M3").Value = -0.8262

$ws.Range("H4").Value = -13.5737
$ws.Range("I4").Value = 168
Set-TextValue $ws.Range("L4") "4.067e+25"
$ws.Range("M4").Value = -1.2263

$ws.Range("H5").Value = -10.6019
$ws.Range("I5").Value = 168
Set-TextValue $ws.Range("L5") "2.359e+17"
$ws.Range("M5").Value = -1.097

$ws.Range("H6").Value = -7.5767
$ws.Range("I6").Value = 168
Set-TextValue $ws.Range("L6") "3.177e+09"
$ws.Range("M6").Value = -0.8264

$ws.Range("H7").Value = -8.192399999999999
$ws.Range("I7").Value = 168
Set-TextValue $ws.Range("L7") "1.05e+11"
$ws.Range("M7").Value = -0.4373

$ws.Range("H8").Value = -4.5678
$ws.Range("I8").Value = 168
Set-TextValue $ws.Range("L8") "1372.761"
$ws.Range("M8").Value = -0.3596

$ws.Range("H9").Value = -0.8544
$ws.Range("I9").Value = 168
$ws.Range("K9").Value = 0.3941
Set-TextValue $ws.Range("L9") "0.123"
$ws.Range("M9").Value = -0.0784

$ws.Range("H10").Value = 1.1264
$ws.Range("I10").Value = 168
$ws.Range("K10").Value = 0.2616
Set-TextValue $ws.Range("L10") "0.16"
$ws.Range("M10").Value = 0.0552

$ws.Range("H11").Value = 5.212
$ws.Range("I11").Value = 168
Set-TextValue $ws.Range("L11") "2.076e+04"
$ws.Range("M11").Value = 0.3334

$ws.Range("H12").Value = 7.6716
$ws.Range("I12").Value = 168
Set-TextValue $ws.Range("L12") "5.4e+09"
$ws.Range("M12").Value = 0.2659

$ws.Range("H13").Value = 0.1835
$ws.Range("I13").Value = 99.96510000000001
$ws.Range("K13").Value = 0.8548
Set-TextValue $ws.Range("L13") "0.19"

$ws.Range("H14").Value = -0.075
$ws.Range("I14").Value = 100.3238
$ws.Range("K14").Value = 0.9404
Set-TextValue $ws.Range("L14") "0.188"
$ws.Range("M14").Value = -0.0114

$ws.Range("H15").Value = 0.0156
$ws.Range("I15").Value = 103.8146
$ws.Range("K15").Value = 0.9875
Set-TextValue $ws.Range("L15") "0.187"
$ws.Range("M15").Value = 0.0023

$ws.Range("H16").Value = -0.3483
$ws.Range("I16").Value = 107.3256
$ws.Range("K16").Value = 0.7282999999999999
Set-TextValue $ws.Range("L16") "0.198"
$ws.Range("M16").Value = -0.0515

$ws.Range("H17").Value = 0.4972
$ws.Range("I17").Value = 108.9776
$ws.Range("K17").Value = 0.62
Set-TextValue $ws.Range("L17") "0.21"
$ws.Range("M17").Value = 0.0731

$ws.Range("H18").Value = 0.5132
$ws.Range("I18").Value = 96.6729
$ws.Range("K18").Value = 0.609
Set-TextValue $ws.Range("L18") "0.211"
$ws.Range("M18").Value = 0.0794

Write-Output "applied stats update"
